# Refresh the coinranking.com crypto price/volume snapshot (GitHub Actions run
# on Sun Aug 20 18:48:17 UTC 2023). Cell text is prefixed with a leading
# apostrophe so Excel stores numeric-looking values (e.g. "1.011") as text,
# matching the original inline-string / General-format cells instead of
# letting Excel auto-convert them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.399.17"
$ws.Range("E2").Value = "'  +0.36%  "
$ws.Range("D3").Value = "'1.693.97"
$ws.Range("E3").Value = "'  +0.25%  "
$ws.Range("D4").Value = "'1.011"
$ws.Range("E4").Value = "'  +0.38%  "
$ws.Range("D5").Value = "'219.26"
$ws.Range("E5").Value = "'  +0.06%  "
$ws.Range("D6").Value = "'0.5483"
$ws.Range("E6").Value = "'  +4.35%  "
$ws.Range("E7").Value = "'  +0.33%  "
$ws.Range("D8").Value = "'0.2736"
$ws.Range("E8").Value = "'  +1.27%  "
$ws.Range("D9").Value = "'0.06470"
$ws.Range("E9").Value = "'  +0.47%  "
$ws.Range("D10").Value = "'22.02"
$ws.Range("E10").Value = "'  -0.18%  "
$ws.Range("D11").Value = "'0.07680"
$ws.Range("E11").Value = "'  +2.76%  "
$ws.Range("D12").Value = "'1.695.86"
$ws.Range("E12").Value = "'  +0.35%  "
$ws.Range("D13").Value = "'4.555"
$ws.Range("E13").Value = "'  -0.13%  "
$ws.Range("D14").Value = "'0.5851"
$ws.Range("E14").Value = "'  -0.17%  "
$ws.Range("D15").Value = "'0.000008372"
$ws.Range("E15").Value = "'  -1.82%  "
$ws.Range("D16").Value = "'65.46"
$ws.Range("E16").Value = "'  +1.27%  "
$ws.Range("D17").Value = "'26.411.42"
$ws.Range("E17").Value = "'  +0.22%  "
$ws.Range("D18").Value = "'4.949"
$ws.Range("E18").Value = "'  -0.23%  "
$ws.Range("D19").Value = "'1.011"
$ws.Range("E19").Value = "'  +0.37%  "
$ws.Range("D20").Value = "'10.97"
$ws.Range("E20").Value = "'  +0.81%  "
$ws.Range("D21").Value = "'192.72"
$ws.Range("E21").Value = "'  +1.44%  "
$ws.Range("E22").Value = "'  +0.52%  "
$ws.Range("E23").Value = "'  +0.38%  "
$ws.Range("D24").Value = "'149.32"
$ws.Range("E24").Value = "'  +3.08%  "
$ws.Range("D25").Value = "'0.1329"
$ws.Range("E25").Value = "'  +7.78%  "
$ws.Range("D26").Value = "'7.915"
$ws.Range("E26").Value = "'  +3.16%  "
$ws.Range("D27").Value = "'15.77"
$ws.Range("E27").Value = "'  -0.60%  "
$ws.Range("D28").Value = "'0.06297"
$ws.Range("E28").Value = "'  -5.50%  "
$ws.Range("D29").Value = "'1.395"
$ws.Range("E29").Value = "'  +3.33%  "
$ws.Range("E30").Value = "'  -0.04%  "
$ws.Range("D31").Value = "'3.605"
$ws.Range("E31").Value = "'  +0.46%  "
$ws.Range("D32").Value = "'3.613"
$ws.Range("E32").Value = "'  +1.23%  "
$ws.Range("D33").Value = "'1.685"
$ws.Range("E33").Value = "'  +0.86%  "
$ws.Range("E34").Value = "'  +1.55%  "
$ws.Range("D35").Value = "'0.6144"
$ws.Range("E35").Value = "'  -1.31%  "
$ws.Range("D36").Value = "'2.411"
$ws.Range("E36").Value = "'  +0.83%  "
$ws.Range("D37").Value = "'2.711"
$ws.Range("E37").Value = "'  +0.54%  "
$ws.Range("D38").Value = "'6.194"
$ws.Range("E38").Value = "'  -2.41%  "
$ws.Range("B39").Value = "'VeChain"
$ws.Range("C39").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01640"
$ws.Range("E39").Value = "'  +0.96%  "
$ws.Range("B40").Value = "'Maker"
$ws.Range("C40").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "'1.119.08"
$ws.Range("E40").Value = "'  +1.10%  "
$ws.Range("D41").Value = "'0.8870"
$ws.Range("D43").Value = "'101.87"
$ws.Range("E43").Value = "'  +0.50%  "
$ws.Range("D44").Value = "'1.843.87"
$ws.Range("E44").Value = "'  +0.33%  "
$ws.Range("E45").Value = "'  -1.11%  "
$ws.Range("E46").Value = "'  +1.06%  "
$ws.Range("D47").Value = "'8.207"
$ws.Range("E47").Value = "'  +0.34%  "
$ws.Range("D48").Value = "'1.007"
$ws.Range("E48").Value = "'  -0.11%  "
$ws.Range("D49").Value = "'0.05285"
$ws.Range("E49").Value = "'  +0.32%  "
$ws.Range("B50").Value = "'Mantle"
$ws.Range("C50").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.4305"
$ws.Range("E50").Value = "'  +0.10%  "
$ws.Range("B51").Value = "'Aptos"
$ws.Range("C51").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "'6.104"
$ws.Range("E51").Value = "'  +0.79%  "
